$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. First paragraph: was empty, becomes the right-justified "Member:" line.
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Member: Jared Orange, Marcos Buzn, Wencong Huang"
$p1.Alignment = 2   # wdAlignParagraphRight

# Re-anchor the lone "_GoBack" bookmark around "Buzn" (Word keeps only one
# bookmark per name in a document, so adding it here also removes it from
# its old location at the end of the "PC:" paragraph).
$bm = $d.Content
$bm.Find.Execute("Buzn", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$d.Bookmarks.Add("_GoBack", $bm)

# ---------------------------------------------------------------------------
# 2. Answer the "PC:" question with the full paragraph of new text, typed in
#    regular (non-bold) weight right after the bold "PC: " label.
# ---------------------------------------------------------------------------
$pcText = "Adequate.  The Program Counter is working with the reset, load, and inc, (the if-else statement in the comment section) total of 3 bits and that gives us 2^3=8 cases. In the test file, it had checked all the cases/combinations of reset, load, and inc.  Therefore, the test is adequate."

$pc = $d.Content
$pc.Find.Execute("PC: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pcStart = $pc.Start
$pcEnd = $pc.End

$pcParagraph = $pcEnd
$p31 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Start -le $pcEnd -and $pcEnd -le $cand.Range.End) {
        $p31 = $cand
        break
    }
}

# Turn off bold for the whole paragraph (text + paragraph mark) so that the
# new answer text -- and the paragraph mark itself -- are not bold, then
# restore bold just on the "PC: " label.
$p31.Range.Bold = 0
$boldRange = $d.Range($pcStart, $pcEnd)
$boldRange.Bold = 1

$insertPoint = $d.Range($pcEnd, $pcEnd)
$insertPoint.InsertAfter($pcText)
$newEnd = $pcEnd + $pcText.Length
$newRange = $d.Range($pcEnd, $newEnd)
$newRange.Bold = 0
